$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the "Node.JS" -> "Node.js" casing typo (unique, case-sensitive
#    match so it cannot collide with the existing "Node.js" spellings).
# ------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Node.JS", $true, $false, $false, $false, $false, $true, 1, $false, "Node.jsX", 2)
# $r now spans the replaced text "Node.jsX" (the trailing "X" is a
# scratch placeholder used below to get a clean, non-paragraph-edge
# insertion point for the bookmark; it is removed again afterwards).

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from its old location (right after the
#    "é" run in the Node.js intro paragraph) to right after the fixed
#    "Node.js" heading run.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$insertPos = $r.End - 1
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the scratch "X" placeholder character again.
$dummyRange = $d.Range($r.End - 1, $r.End)
$dummyRange.Delete()
